# B6-PowerPoint.pptx edit
#
# 1) The three data tables (slides 14, 15, 16) switch from the default
#    "Table_0" table style to the built-in table style
#    {5FEA625B-7680-47D6-A34A-9F1D4B0FDDC4}.
#
# 2) The deck's theme palette (the theme actually applied to the slide
#    master / slides) changes from the "Integral" ("Red Violet") scheme
#    to the default "Office Theme" scheme.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------
$tableStyleId = "{5FEA625B-7680-47D6-A34A-9F1D4B0FDDC4}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($tableStyleId)
        }
    }
}

# --- 2. Swap the slide theme's colour scheme to "Office Theme" ------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (OLE RGB long values,
# i.e. 0x00BBGGRR, matching srgbClr 000000/FFFFFF/44546A/E7E6E6/
# 5B9BD5/ED7D31/A5A5A5/FFC000/4472C4/70AD47/0563C1/954F72).
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$slide1 = $p.Slides.Item(1)
$colorScheme = $slide1.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
